$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric price cells to be stored as text (matches source data which uses
# dotted-thousands / plain decimal strings, not native numbers).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "93.488.88"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.430.04"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "234.34"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "622.10"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +3.87%  "
$ws.Range("D11").Value = "3.426.49"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "43.27"
$ws.Range("E12").Value = "  +7.09%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").Value = "6.34"
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("D15").Value = "93.261.31"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "4.065.92"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "8.26"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "3.426.04"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "18.27"
$ws.Range("E20").Value = "  +7.75%  "
$ws.Range("D21").Value = "11.76"
$ws.Range("E21").Value = "  +6.21%  "
$ws.Range("D22").Value = "3.40"
$ws.Range("E22").Value = "  +8.43%  "
$ws.Range("D23").Value = "502.07"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("D24").Value = "0.479"
$ws.Range("E24").Value = "  +5.19%  "
$ws.Range("D25").Value = "6.66"
$ws.Range("E25").Value = "  +8.26%  "
$ws.Range("D26").Value = "0.0000186"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "95.31"
$ws.Range("E27").Value = "  +6.22%  "
$ws.Range("D28").Value = "12.16"
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").Value = "3.602.68"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "11.45"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "0.174"
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.557"
$ws.Range("E36").Value = "  +4.89%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "29.50"
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("D38").Value = "571.97"
$ws.Range("E38").Value = "  +7.65%  "
$ws.Range("D39").Value = "7.52"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").Value = "1.42"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +1.88%  "
$ws.Range("D43").Value = "0.910"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "23.67"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "3.69"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +5.05%  "
$ws.Range("D48").Value = "5.52"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "53.43"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "8.19"
$ws.Range("E51").Value = "  +4.93%  "
